# This sheet is a weekly price report. Each week two new rows (quality
# "Primera" and "Segunda") are inserted right above the previous week's
# entries, pushing all the older rows down by two. This edit adds this
# week's entries (2023-04-05, serial 45021) as the new rows 436-437.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 436-437; everything from the old row 436 onward
# shifts down to 438 onward (formatting of the row above, e.g. the date
# style on column D, is carried down automatically by Insert()).
$ws.Rows("436:437").Insert()

# Row 436: Primera quality, week of 2023-04-05
$ws.Cells.Item(436, 1).Value = 1
$ws.Cells.Item(436, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(436, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(436, 4).Value = 45021
$ws.Cells.Item(436, 5).Value = 15
$ws.Cells.Item(436, 6).Value = 100112043
$ws.Cells.Item(436, 7).Value = "Pepino ensalada"
$ws.Cells.Item(436, 8).Value = "Sin especificar"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 250
$ws.Cells.Item(436, 11).Value = 5000
$ws.Cells.Item(436, 12).Value = 6000
$ws.Cells.Item(436, 13).Value = 5600
$ws.Cells.Item(436, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(436, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(436, 16).Value = 80
$ws.Cells.Item(436, 17).Value = 70
$ws.Cells.Item(436, 18).Value = "Hortaliza"

# Row 437: Segunda quality, week of 2023-04-05
$ws.Cells.Item(437, 1).Value = 1
$ws.Cells.Item(437, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(437, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(437, 4).Value = 45021
$ws.Cells.Item(437, 5).Value = 15
$ws.Cells.Item(437, 6).Value = 100112043
$ws.Cells.Item(437, 7).Value = "Pepino ensalada"
$ws.Cells.Item(437, 8).Value = "Sin especificar"
$ws.Cells.Item(437, 9).Value = "Segunda"
$ws.Cells.Item(437, 10).Value = 300
$ws.Cells.Item(437, 11).Value = 4000
$ws.Cells.Item(437, 12).Value = 5000
$ws.Cells.Item(437, 13).Value = 4500
$ws.Cells.Item(437, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(437, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(437, 16).Value = 45
$ws.Cells.Item(437, 17).Value = 100
$ws.Cells.Item(437, 18).Value = "Hortaliza"
